$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells keep their literal text (avoid numeric coercion)
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Row 16 and 17 swap (ShibaInu <-> WrappedEther) ---
Set-TextValue $ws.Range("B16") 'WrappedEther'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D16") '2.622.58'
Set-TextValue $ws.Range("E16") '  +0.66%  '

Set-TextValue $ws.Range("B17") 'ShibaInu'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D17") '0.0000133'
Set-TextValue $ws.Range("E17") '  +0.57%  '

# --- Remaining rows: price (D) and volume-1h (E) updates ---
Set-TextValue $ws.Range("D2") '59.504.82'
Set-TextValue $ws.Range("E2") '  +1.06%  '
Set-TextValue $ws.Range("D3") '2.607.75'
Set-TextValue $ws.Range("E3") '  +0.96%  '
Set-TextValue $ws.Range("E4") '  -0.11%  '
Set-TextValue $ws.Range("D5") '537.08'
Set-TextValue $ws.Range("E5") '  +3.08%  '
Set-TextValue $ws.Range("D6") '141.52'
Set-TextValue $ws.Range("E6") '  +2.03%  '
Set-TextValue $ws.Range("E7") '  +0.04%  '
Set-TextValue $ws.Range("E8") '  +0.72%  '
Set-TextValue $ws.Range("E9") '  -1.03%  '
Set-TextValue $ws.Range("E10") '  +1.62%  '
Set-TextValue $ws.Range("E11") '  +1.75%  '
Set-TextValue $ws.Range("E12") '  -0.62%  '
Set-TextValue $ws.Range("D13") '3.068.40'
Set-TextValue $ws.Range("E13") '  +0.92%  '
Set-TextValue $ws.Range("D14") '59.428.71'
Set-TextValue $ws.Range("E14") '  +0.79%  '
Set-TextValue $ws.Range("E15") '  +1.47%  '
Set-TextValue $ws.Range("D18") '341.42'
Set-TextValue $ws.Range("E18") '  +1.39%  '
Set-TextValue $ws.Range("E19") '  +1.81%  '
Set-TextValue $ws.Range("E20") '  +0.30%  '
Set-TextValue $ws.Range("E21") '  -2.05%  '
Set-TextValue $ws.Range("E22") '  -0.02%  '
Set-TextValue $ws.Range("D23") '67.48'
Set-TextValue $ws.Range("E23") '  +2.15%  '
Set-TextValue $ws.Range("E24") '  +1.74%  '
Set-TextValue $ws.Range("E25") '  -1.42%  '
Set-TextValue $ws.Range("E26") '  +0.02%  '
Set-TextValue $ws.Range("E27") '  +3.41%  '
Set-TextValue $ws.Range("D28") '0.0₃0747'
Set-TextValue $ws.Range("E28") '  +3.40%  '
Set-TextValue $ws.Range("E29") '  -0.01%  '
Set-TextValue $ws.Range("E31") '  -1.88%  '
Set-TextValue $ws.Range("E32") '  +1.12%  '
Set-TextValue $ws.Range("D33") '149.57'
Set-TextValue $ws.Range("E33") '  +0.37%  '
Set-TextValue $ws.Range("E34") '  +0.29%  '
Set-TextValue $ws.Range("D35") '1.12'
Set-TextValue $ws.Range("E35") '  -0.10%  '
Set-TextValue $ws.Range("E36") '  +0.47%  '
Set-TextValue $ws.Range("D37") '0.836'
Set-TextValue $ws.Range("E37") '  +3.84%  '
Set-TextValue $ws.Range("D38") '0.828'
Set-TextValue $ws.Range("E38") '  +0.75%  '
Set-TextValue $ws.Range("E39") '  +1.01%  '
Set-TextValue $ws.Range("D40") '1.00'
Set-TextValue $ws.Range("E40") '  +0.03%  '
Set-TextValue $ws.Range("D41") '274.44'
Set-TextValue $ws.Range("E41") '  +1.64%  '
Set-TextValue $ws.Range("E42") '  +2.07%  '
Set-TextValue $ws.Range("E43") '  -0.36%  '
Set-TextValue $ws.Range("E44") '  +0.16%  '
Set-TextValue $ws.Range("E45") '  +1.45%  '
Set-TextValue $ws.Range("D46") '1.948.10'
Set-TextValue $ws.Range("E46") '  -0.78%  '
Set-TextValue $ws.Range("D47") '18.53'
Set-TextValue $ws.Range("E47") '  +3.72%  '
Set-TextValue $ws.Range("E48") '  +1.83%  '
Set-TextValue $ws.Range("D49") '4.51'
Set-TextValue $ws.Range("E49") '  +0.20%  '
Set-TextValue $ws.Range("D50") '110.78'
Set-TextValue $ws.Range("E50") '  -2.29%  '
Set-TextValue $ws.Range("E51") '  +0.46%  '
